# edit.ps1 - Applies the documented changes to the Project System Manual.
#
# Summary of changes:
#   1. "Node is required..." -> "Node/npm is required..." text tweak.
#   2. Swap the two inline image names (image1.png <-> image2.png) on the
#      on-chain / off-chain contract-address screenshots.
#   3. Insert two new bullet paragraphs + a page-break paragraph right after
#      the "Line 7, change const storage_directory ..." bullet (before the
#      "Apart from the above... node corpus_hash_server.js" bullet).
#   4. Insert a new "Please note: ..." bullet right after the
#      "node corpus_hash_server.js" bullet.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simple text replacement
# ---------------------------------------------------------------------
[void]$d.Content.Find.Execute(
    "Node is required to run the provided Corpus Administrator Frontend",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Node/npm is required to run the provided Corpus Administrator Frontend",
    2)

# ---------------------------------------------------------------------
# 2) Swap the image names on the two inline screenshots.
#    Shape 1 (docPr id=22) currently "image1.png" -> "image2.png"
#    Shape 2 (docPr id=21) currently "image2.png" -> "image1.png"
#    Renaming is done by rewriting the shape's own WordOpenXML in place
#    (delete the shape's 1-character range, then re-insert the modified
#    XML at that same, now-collapsed, position).
# ---------------------------------------------------------------------
$shp1 = $d.InlineShapes.Item(1)
$r1 = $shp1.Range
$xml1 = $r1.WordOpenXML
$xml1 = $xml1.Replace('name="image1.png"', 'name="image2.png"')
[void]$r1.Delete()
[void]$r1.InsertXML($xml1)

$shp2 = $d.InlineShapes.Item(2)
$r2 = $shp2.Range
$xml2 = $r2.WordOpenXML
$xml2 = $xml2.Replace('name="image2.png"', 'name="image1.png"')
[void]$r2.Delete()
[void]$r2.InsertXML($xml2)

# ---------------------------------------------------------------------
# Helper: locate a paragraph by unique text, expand the range to cover
# the whole paragraph, and return a Range collapsed to just *before*
# the paragraph's trailing paragraph mark - this is the correct spot to
# InsertXML new sibling <w:p> block(s) right after that paragraph without
# disturbing it or merging into the following paragraph.
# ---------------------------------------------------------------------
function Get-InsertionPointAfterParagraph($doc, $searchText) {
    $fr = $doc.Content
    [void]$fr.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $fr.Expand(4) | Out-Null
    $insertAt = $doc.Range($fr.End - 1, $fr.End - 1)
    return $insertAt
}

# ---------------------------------------------------------------------
# 3) Insert the three new paragraphs after the "Line 7, change const
#    storage_directory ..." bullet.
# ---------------------------------------------------------------------
$pos1 = Get-InsertionPointAfterParagraph $d "Line 7, change const storage_directory = '/upload/corpus/' to:"

$xmlInsert1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="2160" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">This directory must exist and must be accessible from the server that is running the frontend.</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:ind w:left="2160" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">If using the default directory path, </w:t></w:r><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">please note</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> that this directory must be created and </w:t></w:r><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">does not exist</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> in the github repo.</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:left="1440" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:br w:type="page"/></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$pos1.InsertXML($xmlInsert1)

# ---------------------------------------------------------------------
# 4) Insert the new "Please note: ..." paragraph after the
#    "node corpus_hash_server.js" bullet.
# ---------------------------------------------------------------------
$pos2 = Get-InsertionPointAfterParagraph $d "node corpus_hash_server.js"

$xmlInsert2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:left="1440" w:hanging="360"/><w:rPr><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Please note</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">: with this version, any documents uploaded to this version of the system will be stored under the name they were uploaded with. This means that if a document is uploaded that has the same name as an existing document, the existing document will be overwritten.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$pos2.InsertXML($xmlInsert2)

Write-Host "Edits applied successfully."
